$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 111
$ws.Range("B2").Value = "https://insights.blackcoffer.com/ai-and-ml-based-youtube-analytics-and-content-creation-tool-for-optimizing-subscriber-engagement-and-content-strategy/"
$ws.Range("C2").Value = -0.033
$ws.Range("D2").Value = 0.311
$ws.Range("E2").Value = 18.667
$ws.Range("F2").Value = 0.357
$ws.Range("G2").Value = 7.61
$ws.Range("H2").Value = 40
$ws.Range("I2").Value = 112
$ws.Range("J2").Value = 2.188
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 6.857

# Update row 3
$ws.Range("A3").Value = 112
$ws.Range("B3").Value = "https://insights.blackcoffer.com/enhancing-front-end-features-and-functionality-for-improved-user-experience-and-dashboard-accuracy-in-partner-hospital-application/"
$ws.Range("C3").Value = 0.16
$ws.Range("D3").Value = 0.445
$ws.Range("E3").Value = 13.167
$ws.Range("F3").Value = 0.278
$ws.Range("G3").Value = 5.378
$ws.Range("H3").Value = 132
$ws.Range("I3").Value = 474
$ws.Range("J3").Value = 2.093
$ws.Range("K3").Value = 7
$ws.Range("L3").Value = 6.308

# Delete row 4 entirely (shifts nothing up since it's the last row)
$ws.Rows.Item(4).Delete()
